$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits between the
#    "R" run and the "esult..." run of the last paragraph. Deleting it merges
#    the surrounding text back together (the runs stay distinct, matching the
#    target structure) and frees up the name "_GoBack" for re-use later.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Append three new bullet paragraphs after the last ("Result"...) paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.Text = "Visual Aids in analyzing data"
$p1.Range.ListFormat.ListLevelNumber = 1

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.Text = "Experiment set up"
$p2.Range.ListFormat.ListLevelNumber = 2

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3.Range.Text = "Orientation and determining quadrants"
$p3.Range.ListFormat.ListLevelNumber = 2

# 3. Re-create the "_GoBack" bookmark, collapsed right after the text we just
#    typed (before the paragraph mark) in that final new paragraph.
#    A collapsed range placed exactly at "paragraph end - 1" confuses the
#    bookmark engine, so we temporarily type a placeholder character after
#    the text, anchor the bookmark just before it (now a safe mid-paragraph
#    position), and then remove the placeholder again; the bookmark sticks.
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPos = $p3.Range.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()
